$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Förändrad" (C) column for the existing rows (2-399): 45182 -> 45184
$ws.Range("C2:C399").Value = 45184

# --- Touch row 399's height so it gets an explicit custom-height flag, matching the
#     source edit (row 399 picks up ht="15" customHeight="1").
$ws.Rows.Item(399).RowHeight = 15

# --- Append the new record as row 400
$ws.Range("A400").Value = "A 43128-2023"
$ws.Range("B400").Value = 45182
$ws.Range("C400").Value = 45184
$ws.Range("D400").Value = "VÄSTERNORRLANDS LÄN"
$ws.Range("E400").Value = "TIMRÅ"
$ws.Range("F400").Value = "SCA"
$ws.Range("G400").Value = 0.8
$ws.Range("H400:Q400").Value = 0

# Copy the date-number-format from an existing date cell onto the new row's date cells
$ws.Range("C399").Copy()
$ws.Range("B400:C400").PasteSpecial(-4122)

# Copy the wrap-text style used on column R onto the new row's (empty) R cell
$ws.Range("R399").Copy()
$ws.Range("R400").PasteSpecial(-4122)
